# Refresh the crypto price/volume snapshot (Mon Aug 14 09:51:14 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells hold plain text such as "0.9973" or "1.000" that
# would otherwise be auto-coerced into numbers (dropping the trailing zero,
# losing precision, etc.) if written as a bare value, so each is written with
# a leading apostrophe to force literal-text entry, then the style is reset
# back to "Normal" so no stray number-format/quote-prefix style sticks to the
# cell (matching the original, un-styled cells).
$ws.Range("D2").Value = "'29.363.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.847.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.9973"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'240.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.6284"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.07501"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.2903"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'24.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07743"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.846.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'4.992"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.6806"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.00001053"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'82.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'6.191"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'29.406.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'229.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'7.492"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'158.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.1374"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'17.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.06413"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'1.427"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Value = "'4.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'4.097"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Value = "'0.6983"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'2.583"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'1.272.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'2.843"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'6.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.9121"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.9997"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'2.009.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'101.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'66.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'1.737"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Value = "'0.1165"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'9.018"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.3966"
$ws.Range("D50").Style = "Normal"

# Column E ("Volume(1h)") cells are percentage text like "  +0.05%  " -
# never numeric-looking, so a direct assignment is safe.
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  +14.50%  "
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  -18.35%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  -3.40%  "
